$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append two new daily rows (58, 59) for 2025-12-29 (Excel serial date 46020),
# one per station, following the same layout/formatting as the preceding rows.

# Row 58: 四方坪站 -- copy formatting from the row above (row 56, same station)
$ws.Range("A56:F56").Copy()
$ws.Range("A58:F58").PasteSpecial(-4122)
$ws.Range("A58").Value = 46020
$ws.Range("B58").Value = "四方坪站"
$ws.Range("C58").Value = 9427.82
$ws.Range("D58").Value = 8073.52
$ws.Range("E58").Value = 3148.4
$ws.Range("F58").Value = 390

# Row 59: 高岭站 -- copy formatting from the row above (row 57, same station)
$ws.Range("A57:F57").Copy()
$ws.Range("A59:F59").PasteSpecial(-4122)
$ws.Range("A59").Value = 46020
$ws.Range("B59").Value = "高岭站"
$ws.Range("C59").Value = 5636.87
$ws.Range("D59").Value = 4667.8
$ws.Range("E59").Value = 1547.94
$ws.Range("F59").Value = 189

$excel.CutCopyMode = $false

# Move the active-cell selection down past the new data, matching the
# author's post-edit cursor position.
$excel.Goto($ws.Range("I60"))
